# Updated cryptos list on Sat Oct 14 08:34:30 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($Sheet, $Address, $Text)
    $cell = $Sheet.Range($Address)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "26.903.21"
Set-TextValue $ws "E2" "  -0.28%  "
Set-TextValue $ws "D3" "1.549.87"
Set-TextValue $ws "E3" "  -0.46%  "
Set-TextValue $ws "E4" "  -0.32%  "
Set-TextValue $ws "D5" "206.60"
Set-TextValue $ws "E5" "  -0.18%  "
Set-TextValue $ws "E6" "  +0.11%  "
Set-TextValue $ws "E7" "  -0.30%  "
Set-TextValue $ws "D8" "22.22"
Set-TextValue $ws "E8" "  +3.05%  "
Set-TextValue $ws "E9" "  -0.87%  "
Set-TextValue $ws "E10" "  +0.55%  "
Set-TextValue $ws "E11" "  -0.44%  "
Set-TextValue $ws "D12" "1.770.46"
Set-TextValue $ws "D13" "1.549.66"
Set-TextValue $ws "E13" "  -0.34%  "
Set-TextValue $ws "E14" "  +0.78%  "
Set-TextValue $ws "E15" "  +0.79%  "
Set-TextValue $ws "D16" "26.906.29"
Set-TextValue $ws "E16" "  -0.23%  "
Set-TextValue $ws "D17" "61.74"
Set-TextValue $ws "E17" "  -0.07%  "
Set-TextValue $ws "D18" "217.69"
Set-TextValue $ws "E18" "  +1.31%  "
Set-TextValue $ws "E19" "  +1.70%  "
Set-TextValue $ws "E20" "  -0.04%  "
Set-TextValue $ws "E21" "  -0.30%  "
Set-TextValue $ws "E22" "  +0.13%  "
Set-TextValue $ws "D23" "9.23"
Set-TextValue $ws "E23" "  +0.12%  "
Set-TextValue $ws "E24" "  -1.00%  "
Set-TextValue $ws "D25" "154.30"
Set-TextValue $ws "E25" "  +0.55%  "
Set-TextValue $ws "E26" "  -0.35%  "
Set-TextValue $ws "E27" "  +0.12%  "
Set-TextValue $ws "E28" "  +0.74%  "
Set-TextValue $ws "E30" "  +1.56%  "
Set-TextValue $ws "E31" "  -0.89%  "
Set-TextValue $ws "E32" "  -0.22%  "
Set-TextValue $ws "D33" "1.416.61"
Set-TextValue $ws "E33" "  +3.12%  "
Set-TextValue $ws "E34" "  +3.91%  "
Set-TextValue $ws "D35" "1.58"
Set-TextValue $ws "E35" "  +2.07%  "
Set-TextValue $ws "D36" "0.965"
Set-TextValue $ws "E36" "  -0.31%  "
Set-TextValue $ws "E37" "  -0.08%  "
Set-TextValue $ws "E38" "  +0.03%  "
Set-TextValue $ws "D39" "0.527"
Set-TextValue $ws "E39" "  +1.28%  "
Set-TextValue $ws "E40" "  -0.10%  "
Set-TextValue $ws "B41" "FraxShare"
Set-TextValue $ws "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D41" "5.71"
Set-TextValue $ws "E41" "  +3.68%  "
Set-TextValue $ws "B42" "PaxDollar"
Set-TextValue $ws "C42" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws "D42" "1.00"
Set-TextValue $ws "E42" "  -0.31%  "
Set-TextValue $ws "E43" "  +4.20%  "
Set-TextValue $ws "E44" "  +1.96%  "
Set-TextValue $ws "D45" "64.49"
Set-TextValue $ws "E45" "  +1.05%  "
Set-TextValue $ws "E46" "  -0.28%  "
Set-TextValue $ws "D47" "1.684.30"
Set-TextValue $ws "D48" "87.54"
Set-TextValue $ws "E48" "  +1.39%  "
Set-TextValue $ws "E49" "  +4.35%  "
Set-TextValue $ws "D50" "0.0518"
Set-TextValue $ws "E50" "  +1.86%  "
Set-TextValue $ws "D51" "0.0959"
Set-TextValue $ws "E51" "  +0.29%  "
